# Actualización automática de tasas-transfi.xlsx

$wb = $excel.ActiveWorkbook

# --- Hoja1: actualizar el texto de la conversión del día ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$texto = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 8.52 = 35750.1 pesos`n✅ 35750.1 pesos = 8.52 = 973.1 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $texto

# --- tasas: actualizar tasas N10/O10 y N12/O12 ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 117.398
$wsTasas.Range("O10").Value = 4196.99

$wsTasas.Range("N12").Value = 4198
$wsTasas.Range("O12").Value = 114.267
